{"js": "// Mapping of old text -> new text for every text run in the document\n// (the title date plus every multiplication expression in the table).\nconst replacements = [\n  [\"2025-12-28 Sunday\", \"2025-12-29 Monday\"],\n  [\"23\u00d733=759\", \"88\u00d793=8184\"],\n  [\"82\u00d790=7380\", \"42\u00d757=2394\"],\n  [\"72\u00d721=1512\", \"32\u00d793=2976\"],\n  [\"96\u00d779=7584\", \"13\u00d752=676\"],\n  [\"82\u00d785=6970\", \"37\u00d737=1369\"],\n  [\"63\u00d785=5355\", \"73\u00d727=1971\"],\n  [\"58\u00d757=3306\", \"22\u00d758=1276\"],\n  [\"65\u00d720=1300\", \"83\u00d781=6723\"],\n  [\"54\u00d714=756\", \"90\u00d785=7650\"],\n  [\"77\u00d762=4774\", \"40\u00d715=600\"],\n  [\"11\u00d755=605\", \"26\u00d787=2262\"],\n  [\"87\u00d735=3045\", \"43\u00d716=688\"],\n  [\"19\u00d771=1349\", \"41\u00d719=779\"],\n  [\"28\u00d746=1288\", \"28\u00d791=2548\"],\n  [\"44\u00d725=1100\", \"52\u00d795=4940\"],\n  [\"94\u00d726=2444\", \"30\u00d711=330\"],\n  [\"38\u00d768=2584\", \"93\u00d795=8835\"],\n  [\"72\u00d723=1656\", \"29\u00d751=1479\"],\n  [\"35\u00d768=2380\", \"60\u00d768=4080\"],\n  [\"34\u00d748=1632\", \"30\u00d754=1620\"],\n  [\"19\u00d750=950\", \"36\u00d735=1260\"],\n  [\"83\u00d766=5478\", \"24\u00d751=1224\"],\n  [\"99\u00d733=3267\", \"82\u00d746=3772\"],\n  [\"71\u00d773=5183\", \"23\u00d727=621\"],\n  [\"61\u00d797=5917\", \"14\u00d736=504\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Mapping of old text -> new text for every text run in the document\n# (the title date plus every multiplication expression in the table).\n$replacements = @(\n    @(\"2025-12-28 Sunday\", \"2025-12-29 Monday\"),\n    @(\"23\u00d733=759\", \"88\u00d793=8184\"),\n    @(\"82\u00d790=7380\", \"42\u00d757=2394\"),\n    @(\"72\u00d721=1512\", \"32\u00d793=2976\"),\n    @(\"96\u00d779=7584\", \"13\u00d752=676\"),\n    @(\"82\u00d785=6970\", \"37\u00d737=1369\"),\n    @(\"63\u00d785=5355\", \"73\u00d727=1971\"),\n    @(\"58\u00d757=3306\", \"22\u00d758=1276\"),\n    @(\"65\u00d720=1300\", \"83\u00d781=6723\"),\n    @(\"54\u00d714=756\", \"90\u00d785=7650\"),\n    @(\"77\u00d762=4774\", \"40\u00d715=600\"),\n    @(\"11\u00d755=605\", \"26\u00d787=2262\"),\n    @(\"87\u00d735=3045\", \"43\u00d716=688\"),\n    @(\"19\u00d771=1349\", \"41\u00d719=779\"),\n    @(\"28\u00d746=1288\", \"28\u00d791=2548\"),\n    @(\"44\u00d725=1100\", \"52\u00d795=4940\"),\n    @(\"94\u00d726=2444\", \"30\u00d711=330\"),\n    @(\"38\u00d768=2584\", \"93\u00d795=8835\"),\n    @(\"72\u00d723=1656\", \"29\u00d751=1479\"),\n    @(\"35\u00d768=2380\", \"60\u00d768=4080\"),\n    @(\"34\u00d748=1632\", \"30\u00d754=1620\"),\n    @(\"19\u00d750=950\", \"36\u00d735=1260\"),\n    @(\"83\u00d766=5478\", \"24\u00d751=1224\"),\n    @(\"99\u00d733=3267\", \"82\u00d746=3772\"),\n    @(\"71\u00d773=5183\", \"23\u00d727=621\"),\n    @(\"61\u00d797=5917\", \"14\u00d736=504\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
